$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Heath Mellen"
$ws.Range("B2").Value = "HM"

$ws.Range("A3").Value = "James Aaron"
$ws.Range("B3").Value = "JA"

$ws.Range("A4").Value = "Laura Jackson"
$ws.Range("B4").Value = "LJ"

$ws.Range("A5").Value = "Alfred Ludwig"
$ws.Range("B5").Value = "AL"

$ws.Range("A6").Value = "David Tarver"
$ws.Range("B6").Value = "DT"

$ws.Range("A7").Value = "Robyn George"
$ws.Range("B7").Value = "RG"

$ws.Range("A8").Value = "Randy Byerly"
$ws.Range("B8").Value = "RB"
